$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3915976666666667
$ws.Range("H2").Value = 1.174793
$ws.Range("I2").Value = 0.02606065131430495
$ws.Range("J2").Value = 0.02606065131430495
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 5.608415181877668
$ws.Range("R2").Value = 50.47573663689901
$ws.Range("S2").Value = 0.007686769384356498
$ws.Range("T2").Value = 0.007686769384356498
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3915976666666667
$ws.Range("H3").Value = 1.174793
$ws.Range("I3").Value = 0.02606065131430495
$ws.Range("J3").Value = 0.02606065131430495
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 10.60609738400567
$ws.Range("R3").Value = 95.45487645605103
$ws.Range("S3").Value = 0.0145364817002695
$ws.Range("T3").Value = 0.0145364817002695
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3915976666666667
$ws.Range("H4").Value = 1.174793
$ws.Range("I4").Value = 0.02606065131430495
$ws.Range("J4").Value = 0.02606065131430495
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 2.799841211689222
$ws.Range("R4").Value = 25.198570905203
$ws.Range("S4").Value = 0.003837400229678959
$ws.Range("T4").Value = 0.003837400229678958
$ws.Range("I5").Value = 0.4187506438669658
$ws.Range("J5").Value = 0.4187506438669658
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 90.117758000753
$ws.Range("R5").Value = 811.059822006777
$ws.Range("S5").Value = 0.1235133991908065
$ws.Range("T5").Value = 0.1235133991908065
$ws.Range("I6").Value = 0.4187506438669658
$ws.Range("J6").Value = 0.4187506438669658
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("S6").Value = 0.2335767052839126
$ws.Range("T6").Value = 0.2335767052839126
$ws.Range("I7").Value = 0.4187506438669658
$ws.Range("J7").Value = 0.4187506438669658
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 44.98871866170767
$ws.Range("R7").Value = 404.898467955369
$ws.Range("S7").Value = 0.0616605393922467
$ws.Range("T7").Value = 0.06166053939224669
$ws.Range("G8").Value = 8.342485333333334
$ws.Range("H8").Value = 25.027456
$ws.Range("I8").Value = 0.5551887048187292
$ws.Range("J8").Value = 0.5551887048187292
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 119.4800821882453
$ws.Range("R8").Value = 1075.320739694208
$ws.Range("S8").Value = 0.1637567491031436
$ws.Range("T8").Value = 0.1637567491031436
$ws.Range("G9").Value = 8.342485333333334
$ws.Range("H9").Value = 25.027456
$ws.Range("I9").Value = 0.5551887048187292
$ws.Range("J9").Value = 0.5551887048187292
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 225.9492826480214
$ws.Range("R9").Value = 2033.543543832192
$ws.Range("S9").Value = 0.3096810724513169
$ws.Range("T9").Value = 0.3096810724513169
$ws.Range("G10").Value = 8.342485333333334
$ws.Range("H10").Value = 25.027456
$ws.Range("I10").Value = 0.5551887048187292
$ws.Range("J10").Value = 0.5551887048187292
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 59.64702099224178
$ws.Range("R10").Value = 536.823188930176
$ws.Range("S10").Value = 0.0817508832642687
$ws.Range("T10").Value = 0.08175088326426869
